$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 10.58
$ws.Range("E3").Value = 10.38
$ws.Range("F3").Value = 9.75
$ws.Range("G3").Value = 11.9

$ws.Range("C4").Value = 9.42
$ws.Range("E4").Value = 10.17
$ws.Range("F4").Value = 10.16

$ws.Range("C5").Value = 9.619999999999999
$ws.Range("D5").Value = 9.83
$ws.Range("F5").Value = 10.31
$ws.Range("G5").Value = 9.210000000000001
$ws.Range("H5").Value = 8.43

$ws.Range("C6").Value = 10.25
$ws.Range("D6").Value = 9.84
$ws.Range("E6").Value = 9.69
$ws.Range("G6").Value = 10.35
$ws.Range("H6").Value = 11.39
$ws.Range("J6").Value = 8.08

$ws.Range("C7").Value = 8.1
$ws.Range("E7").Value = 10.79
$ws.Range("F7").Value = 9.65

$ws.Range("E8").Value = 11.57
$ws.Range("F8").Value = 8.609999999999999

$ws.Range("F10").Value = 11.92
